$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "22.392.02"
$ws.Cells.Item(3, 4).Value = "1.568.87"
$ws.Cells.Item(3, 5).Value = "  -4.70%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 5).Value = "  +0.01%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "291.19"
$ws.Cells.Item(6, 5).Value = "  -2.67%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3664"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "49.44"
$ws.Cells.Item(8, 5).Value = "  -1.21%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.3384"
$ws.Cells.Item(9, 5).Value = "  -4.23%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.173"
$ws.Cells.Item(10, 5).Value = "  -3.34%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07590"
$ws.Cells.Item(12, 5).Value = "  -0.02%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "21.20"
$ws.Cells.Item(13, 5).Value = "  -4.04%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.057"
$ws.Cells.Item(14, 5).Value = "  -5.23%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "6.899"
$ws.Cells.Item(15, 5).Value = "  -5.86%  "
$ws.Cells.Item(16, 5).Value = "  -5.14%  "
$ws.Cells.Item(17, 4).Value = "1.570.69"
$ws.Cells.Item(17, 5).Value = "  -4.83%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "89.06"
$ws.Cells.Item(18, 5).Value = "  -8.15%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06756"
$ws.Cells.Item(19, 5).Value = "  -3.00%  "
$ws.Cells.Item(20, 5).Value = "  +0.05%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.236"
$ws.Cells.Item(21, 5).Value = "  -7.54%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.5311"
$ws.Cells.Item(22, 5).Value = "  -8.33%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "16.44"
$ws.Cells.Item(23, 5).Value = "  -5.51%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "12.07"
$ws.Cells.Item(24, 5).Value = "  -2.82%  "
$ws.Cells.Item(25, 4).Value = "22.398.68"
$ws.Cells.Item(25, 5).Value = "  -4.56%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.401"
$ws.Cells.Item(26, 5).Value = "  -3.99%  "
$ws.Cells.Item(27, 5).Value = "  +4.06%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "19.89"
$ws.Cells.Item(28, 5).Value = "  -4.67%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "144.80"
$ws.Cells.Item(29, 5).Value = "  -4.93%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.965"
$ws.Cells.Item(30, 5).Value = "  -4.47%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "125.34"
$ws.Cells.Item(31, 5).Value = "  -5.59%  "
$ws.Cells.Item(32, 4).Value = "1.744.44"
$ws.Cells.Item(32, 5).Value = "  -4.74%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.038"
$ws.Cells.Item(33, 5).Value = "  +4.74%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.264"
$ws.Cells.Item(34, 5).Value = "  -9.82%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.985"
$ws.Cells.Item(35, 5).Value = "  -7.70%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "10.35"
$ws.Cells.Item(36, 5).Value = "  -9.50%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.02563"
$ws.Cells.Item(37, 5).Value = "  -5.46%  "
$ws.Cells.Item(38, 5).Value = "  -3.21%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.2307"
$ws.Cells.Item(39, 5).Value = "  -5.33%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.06538"
$ws.Cells.Item(40, 5).Value = "  -3.77%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.520"
$ws.Cells.Item(41, 5).Value = "  -6.97%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "11.87"
$ws.Cells.Item(42, 5).Value = "  -8.07%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.252"
$ws.Cells.Item(43, 5).Value = "  -3.27%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.6392"
$ws.Cells.Item(44, 5).Value = "  -7.09%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "14.37"
$ws.Cells.Item(45, 5).Value = "  -8.66%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.000"
$ws.Cells.Item(46, 5).Value = "  -0.03%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.6010"
$ws.Cells.Item(47, 5).Value = "  -5.42%  "
$ws.Cells.Item(48, 5).Value = "  -3.17%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.134"
$ws.Cells.Item(49, 5).Value = "  -5.24%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "123.14"
$ws.Cells.Item(50, 5).Value = "  -3.41%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.213"
$ws.Cells.Item(51, 5).Value = "  +2.80%  "
